$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.108.63'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '1.650.73'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.17'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5229'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.79%  '
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06285'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.48'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07794'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.479'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('D13').Value = '1.655.27'
$ws.Range('E13').Value = '  -1.02%  '
$ws.Range('D14').Value = '1.877.74'
$ws.Range('E14').Value = '  -0.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5529'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '0.0₅8002'
$ws.Range('E16').Value = '  -2.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.71'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').Value = '26.090.72'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.619'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '194.25'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.07'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.945'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.14'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1202'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.155'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.478'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05701'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.268'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.476'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.337'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.585'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.57%  '
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9484'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.413'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5657'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01589'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.924'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.51%  '
$ws.Range('D41').Value = '1.060.94'
$ws.Range('E41').Value = '  +1.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.005'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8418'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.34'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.64%  '
$ws.Range('D45').Value = '1.789.05'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.37'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05398'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.57%  '
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4398'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.972'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.68%  '
